# Fix uom and offer not displaying in joi
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Site PR field (I8): "PR10011234" -> "PR100112341"
$ws.Range("I8").Value = "PR100112341"

# Department field (I7): " department11" -> "  IT-Bacolod"
$ws.Range("I7").Value = "  IT-Bacolod"

# Urgency No. field (I10): empty -> 0
$ws.Range("I10").Value = 0

# Update selection to I7:K7 with active cell I7
$ws.Range("I7:K7").Select()
